# Adding auto-match for invoices and receipts
# Updates the "receipts" sheet: re-labels payment modes, adds a
# "reference" column, bumps one receipt amount, appends a new receipt
# row, and makes "receipts" the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("receipts")

# --- Header: new "reference" column ---
$ws.Range("F1").Value = "reference"

# --- Row 2: capitalize existing payment mode ---
$ws.Range("E2").Value = "Cash"

# --- Row 3: switch payment mode to UPI, add a reference, bump amount ---
$ws.Range("D3").Value = 2600
$ws.Range("E3").Value = "UPI"
$ws.Range("F3").Value = "abc123"

# --- Row 4: brand new receipt ---
$ws.Range("A4").Value = "REC-2025-003"

# copy the date formatting from the row above so the new date cell
# keeps the same built-in date style (s="1") instead of creating a new one
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B4").Value = 45905

$ws.Range("C4").Value = "Spencer Bradford"
$ws.Range("D4").Value = 700
$ws.Range("E4").Value = "Cash"

# --- Make "receipts" the active sheet with D5 selected ---
$ws.Activate()
$ws.Range("D5").Select()
